# Test-Scenarios.xlsx edit: "TC for paymnet added"
#
# Adds two new Payment test-case rows (TS_PAY_06 / TS_PAY_07) and a new
# "different payment method" scenario, and tweaks the selection / active
# sheet state to match the authored commit.

$wb = $excel.ActiveWorkbook

$paySheet = $wb.Worksheets.Item("Payment")

# --- Write new content in the same order the shared-string table grew in
#     the real edit (IDs first, then requirement codes, then the new
#     descriptions) so newly-created shared strings land in the expected
#     slots. ---

# New scenario / requirement IDs for the two appended rows.
$paySheet.Range("A11").Value = "TS_PAY_06"
$paySheet.Range("A12").Value = "TS_PAY_07"
$paySheet.Range("B11").Value = "FR-PAY-06"
$paySheet.Range("B12").Value = "FR-PAY-07"

# A new scenario was inserted ahead of the existing payment-method
# scenarios, so the "Credit/Debit Card" row now carries a fresh
# description while the old descriptions shift down one row.
$paySheet.Range("C6").Value = "Verify different payment method avialable"
$paySheet.Range("C11").Value = "Verify  payment Success"
$paySheet.Range("C12").Value = "Verify secure payment "

# Existing descriptions shift down into the row below their old slot.
$paySheet.Range("C7").Value = "Verify payment using Credit/Debit Card"
$paySheet.Range("C8").Value = "Verify payment using UPI"
$paySheet.Range("C9").Value = "Verify Cash on Delivery option"
$paySheet.Range("C10").Value = "Verify payment failure handling"

# --- Selection / active-sheet bookkeeping, matching the saved view state ---

$cartSheet = $wb.Worksheets.Item("Add to Cart")
[void]$cartSheet.Activate()
[void]$cartSheet.Range("A11").Select()

$checkoutSheet = $wb.Worksheets.Item("Checkout & Address Management")
[void]$checkoutSheet.Activate()
[void]$checkoutSheet.Range("B17").Select()

[void]$paySheet.Activate()
[void]$paySheet.Range("C15").Select()
